$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.169627547264099
$ws.Range("B1").Value = 2.440605640411377
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.362940549850464
$ws.Range("E1").Value = 1.237362027168274
